# Fixed swapped size/express bug
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

# The "Size" labels for the scale rows were showing the wrong text
# (the express/"X" column value had been duplicated into the size column).
# Correct the two scale-cable size labels.
$ws.Range("A13").Value = "SCALE 1"
$ws.Range("A14").Value = "SCALE 2"

# Re-select A14, matching the cursor position left behind by the fix.
$ws.Range("A14").Select()
